$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H), matching the header styling already used by
# the other header cells (B1:G1) - copy G1's format onto H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save value for the single data row.
$ws.Range("H2").Value = 1
